# "add some profiler points, Krispi level up." -- append six new backlog
# rows to the "Бэклог задач" sheet (sheet2), matching new tasks/timestamps
# picked up from the profiler pass, then refresh the sheet's view state,
# the (now much longer) column B width, and the page setup.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- new backlog rows (49-54): description in col B, timestamp in col C ---
$rows = @(
    @{ Row = 49; Text = "Постинг нша стену результатов игры."; Date = 42037.756944444445 },
    @{ Row = 50; Text = "Оптимизация: индексы для users.socNetTypeId и users.socNetUserId"; Date = 42038.629861111112 },
    @{ Row = 51; Text = "Оптимизация: кэширование для запросов к социальной сети! В первую очередь это нужно для тестирования, т.к. мы просто замучаем АПИ соц сети,  не хорошо это, вот."; Date = 42038.629861111112 },
    @{ Row = 52; Text = "Рейтинг пытается отобразить английские буквы, failed"; Date = 42038.938888888886 },
    @{ Row = 53; Text = "ElementGraphicsText если нет картинки символа - заменять на текст."; Date = 42038.950694444444 },
    @{ Row = 54; Text = "При  закрытие клиента, похоже, что его игры не закрываются, а остаются в статусе 2!, если это так, соответствено добавить на это тест."; Date = 42039.758333333331 }
)

foreach ($r in $rows) {
    $ws.Range("B" + $r.Row).Value = $r.Text
    $cell = $ws.Range("C" + $r.Row)
    $cell.Value = $r.Date
    $cell.NumberFormat = "m/d/yy h:mm"
}

# --- column B grew wide enough to need a resize (bestFit to ~176 chars) ---
$ws.Columns.Item(2).ColumnWidth = 175.15

# --- view: last row is now the active cell, scrolled further down ---
$ws.Activate()
$ws.Range("B54").Select() | Out-Null

# --- page setup for the sheet ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
